# Fruta / hortaliza, semanal
# Inserts a new weekly price-report row for Chirimoya (Vega Modelo de Temuco)
# above the current row 189, shifting the existing rows 189-227 down to
# 190-228 (exactly like dragging a new row into the middle of the table in
# Excel), then fills the newly-inserted row with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new blank row at 189; everything that used to live on row
# 189 (and below) moves down by one row, and the sheet's used range grows
# from A1:T227 to A1:T228 automatically.
$ws.Rows(189).Insert()

# Populate the newly inserted row 189 with the new weekly record.
$ws.Range("A189").Value = 10
$ws.Range("B189").Value = "Vega Modelo de Temuco"
$ws.Range("C189").Value = "La Araucanía"
$ws.Range("D189").Value = 45209
$ws.Range("E189").Value = 9
$ws.Range("F189").Value = "Fruta"
$ws.Range("G189").Value = 100107
$ws.Range("H189").Value = "Otros"
$ws.Range("I189").Value = 100107002
$ws.Range("J189").Value = "Chirimoya"
$ws.Range("K189").Value = "Cultivar IV Región"
$ws.Range("L189").Value = "Primera"
$ws.Range("M189").Value = 155
$ws.Range("N189").Value = 2500
$ws.Range("O189").Value = 2500
$ws.Range("P189").Value = 2500
$ws.Range("Q189").Value = "`$/kilo (en caja de 15 kilos)"
$ws.Range("R189").Value = "Provincia del Elquí"
$ws.Range("S189").Value = 2500
$ws.Range("T189").Value = 1
